# SplashScreen.pptx maintenance edit
#  - bump the footer "date" field shown on the slide master + every slide
#    layout from 3/24/2021 -> 3/26/2021
#  - fix the title on the title slide: "Data Lotto Manager" -> "Lotto Data Manager"

$p = $ppt.ActivePresentation

$oldDate = "3/24/2021"
$newDate = "3/26/2021"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shape) {
    if (-not $shape.HasTextFrame) { return }
    try {
        $phType = $shape.PlaceholderFormat.Type
    } catch {
        return
    }
    if ($phType -eq $ppPlaceholderDate) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# 1) Slide master footer date field
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DatePlaceholder $master.Shapes.Item($i)
}

# 2) Every slide layout's footer date field
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        Update-DatePlaceholder $layout.Shapes.Item($si)
    }
}

# 3) Title-slide headline text
$slide = $p.Slides.Item(1)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $shape = $slide.Shapes.Item($si)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -eq "Data Lotto Manager") {
            $shape.TextFrame.TextRange.Text = "Lotto Data Manager"
        }
    }
}
